$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy style (bold, border, centered/top alignment) from the last existing
# data row (A10) onto the new year-label cells A11 and A12.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# Row 11 (2021年)
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 34427.66
$ws.Range("C11").Value = 9115.120000000001
$ws.Range("D11").Value = 1951.23
$ws.Range("E11").Value = 3.93
$ws.Range("F11").Value = 20084.05
$ws.Range("G11").Value = 41273.64
$ws.Range("H11").Value = 4497.34
$ws.Range("I11").Value = 25824.98
$ws.Range("J11").Value = 4002.75
$ws.Range("K11").Value = 4240.34
$ws.Range("L11").Value = 2749.21
$ws.Range("M11").Value = 1639.92
$ws.Range("N11").Value = 6269.25
$ws.Range("O11").Value = 23263.22
$ws.Range("P11").Value = 2158.69
$ws.Range("Q11").Value = 3389.01
$ws.Range("R11").Value = 16276.92
$ws.Range("S11").Value = 7189.01
$ws.Range("T11").Value = 53349.38
$ws.Range("U11").Value = 8390.049999999999
$ws.Range("V11").Value = 31288.97
$ws.Range("W11").Value = 5417.93
$ws.Range("X11").Value = 38743.27
$ws.Range("Y11").Value = 60255.22
$ws.Range("Z11").Value = 4118.86
$ws.Range("AA11").Value = 19328.93
$ws.Range("AB11").Value = 3507.81
$ws.Range("AC11").Value = 12484.67
$ws.Range("AD11").Value = 7124.94
$ws.Range("AE11").Value = 751807.17
$ws.Range("AF11").Value = 98880.62
$ws.Range("AG11").Value = 38805.02
$ws.Range("AH11").Value = 8271.33
$ws.Range("AI11").Value = 12024.29
$ws.Range("AJ11").Value = 1451.74
$ws.Range("AK11").Value = 27164.37
$ws.Range("AL11").Value = 21887.48
$ws.Range("AM11").Value = 41993.29
$ws.Range("AN11").Value = 2948.2
$ws.Range("AO11").Value = 10349.95
$ws.Range("AP11").Value = 31749.28
$ws.Range("AQ11").Value = 3915.3

# Row 12 (2022年)
$ws.Range("A12").Value = "2022年"
$ws.Range("B12").Value = 38620.3
$ws.Range("C12").Value = 9987.5
$ws.Range("D12").Value = 1087.6
$ws.Range("E12").Value = 3.4
$ws.Range("F12").Value = 22460.5
$ws.Range("G12").Value = 45399.7
$ws.Range("H12").Value = 4820.9
$ws.Range("I12").Value = 28080
$ws.Range("J12").Value = 4203.5
$ws.Range("K12").Value = 4351.3
$ws.Range("L12").Value = 3485.6
$ws.Range("M12").Value = 1772.6
$ws.Range("N12").Value = 6584.1
$ws.Range("O12").Value = 26975.3
$ws.Range("P12").Value = 2286.9
$ws.Range("Q12").Value = 3660.1
$ws.Range("R12").Value = 17299.5
$ws.Range("S12").Value = 8142.1
$ws.Range("T12").Value = 59908.8
$ws.Range("U12").Value = 7351.4
$ws.Range("V12").Value = 34494.1
$ws.Range("W12").Value = 5988.8
$ws.Range("X12").Value = 42406.7
$ws.Range("Y12").Value = 73878
$ws.Range("Z12").Value = 4171
$ws.Range("AA12").Value = 19953.9
$ws.Range("AB12").Value = 4054.8
$ws.Range("AC12").Value = 12888.9
$ws.Range("AD12").Value = 7085
$ws.Range("AE12").Value = 807645.9
$ws.Range("AF12").Value = 107487.3
$ws.Range("AG12").Value = 39921.3
$ws.Range("AH12").Value = 8401
$ws.Range("AI12").Value = 13201
$ws.Range("AJ12").Value = 1578.6
$ws.Range("AK12").Value = 25519.7
$ws.Range("AL12").Value = 13391.7
$ws.Range("AM12").Value = 45573.7
$ws.Range("AN12").Value = 3351.3
$ws.Range("AO12").Value = 11891.7
$ws.Range("AP12").Value = 32162.7
$ws.Range("AQ12").Value = 3763.8
